$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (MH007): price 39999 -> 19000, status Hết hàng -> Còn hàng
$ws.Range("C8").Value = 19000.0
$ws.Range("D8").Value = "Còn hàng"

# Row 10 (MH009): price 23333 -> 12000, status Hết hàng -> Còn hàng
$ws.Range("C10").Value = 12000.0
$ws.Range("D10").Value = "Còn hàng"

# Row 11 (MH010): price 12222 -> 12000, status Hết hàng -> Còn hàng
$ws.Range("C11").Value = 12000.0
$ws.Range("D11").Value = "Còn hàng"

# Row 12 (MH011 / Nhi Nhi) is removed entirely - delete the whole row
$ws.Rows.Item(12).Delete()
